$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Toplam" (Total) column - column F - is dropped from the export.
# Deleting it shifts every following column one slot to the left, which is
# exactly the before/after realignment seen in the diff:
#   F: Toplam -> Muhasebe Kodu
#   G: Muhasebe Kodu -> Belge No
#   H: Belge No -> Cek no
#   I: Cek no -> Para birim(i)
#   J: Para birim -> Para Birimi Tutari
#   K: Para birimi tutari -> (now empty, dimension shrinks by one column)
$ws.Columns.Item(6).Delete()

# Escape-char / typo fixes to the two currency-related headers, now sitting
# in columns I and J after the shift above.
$ws.Cells.Item(1, 9).Value = "Para birimi"
$ws.Cells.Item(1, 10).Value = "Para Birimi Tutarı"

# The active selection moves to the new "Para birimi" header cell's column.
$ws.Range("J2").Select()

# Nudge the workbook's far-right (unused) column width slightly, matching
# the new <col max="1024" min="1024" .../> entry that appears after the edit.
$ws.Columns.Item(1024).ColumnWidth = 10.6
